$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Registrar pago de alumno: la tarea de vincular pago de alumno con grupo
# ya fue completada (Hecho) y se consumieron 3 horas el primer dia.
$ws.Range("F11").Value = "Hecho"
$ws.Range("H11").Value = 3

$ws.Range("F11").Select()
